$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (paragraph 1). It is built with InsertXML so we get the exact same
#    run layout as the target: a leading empty run, a bold "Meta
#    description" run, and a plain run with the rest of the sentence.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$newPara = $titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr>
                <w:b/>
              </w:rPr>
              <w:t>Meta description</w:t>
            </w:r>
            <w:r>
              <w:t>: Discover Dog Town Deal, a charming slot game with excellent graphics, a jazz soundtrack, and low volatility gameplay. Play for free and find out if the payout potential is worth it.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$null = $metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph that sits at the very end
#    of the document (right before the final italic paragraph). Walk the
#    paragraphs from the end so this is robust no matter how the earlier
#    insert shifted paragraph indices.
# ---------------------------------------------------------------------------
$dupTitleText = "Play Dog Town Deal for Free - Review of Adorable Slot Game"
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd() -eq $dupTitleText) {
        $candidate.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph - the one that still
#    holds the old "Discover Dog Town Deal..." description - with the new
#    image prompt copy, keeping its existing italic run formatting intact.
# ---------------------------------------------------------------------------
$oldDescriptionText = "Discover Dog Town Deal, a charming slot game with excellent graphics, a jazz soundtrack, and low volatility gameplay. Play for free and find out if the payout potential is worth it."
$newDescriptionText = 'Create a feature image fitting the game "Dog Town Deal" with the following specifications: Design Prompt: Create a cartoon-style image featuring a happy Maya warrior with glasses. The Maya warrior should be surrounded by lovable canines playing the role of gamblers, with chips stamped with a paw print and a stylized bone Wild symbol. The background should be set in a clandestine den with elegantly dressed dogs standing beside the game table. The image should be well-detailed and full of life, much like the game itself. Ensure that the image focuses on the theme, giving players an idea of what to expect from the game.'

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd() -eq $oldDescriptionText) {
        $textRange = $d.Range($candidate.Range.Start, $candidate.Range.End - 1)
        $textRange.Text = $newDescriptionText
        break
    }
}
